$d = $word.ActiveDocument

# Paragraph "{{LETRA_G}} ..." - the phrase "QUE EL PROMITENTE VENDEDOR LE ENTREGUE"
# is replaced with "QUE "{{SEXO_1}} PROMITENTE {{SEXO_2}}" LE ENTREGUE"
$d.Content.Find.Execute(
    "QUE EL PROMITENTE VENDEDOR LE ENTREGUE",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "QUE " + [char]8220 + "{{SEXO_1}} PROMITENTE {{SEXO_2}}" + [char]8221 + " LE ENTREGUE",
    2
)
